$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1635.9286
$ws.Range("I100").Value = 1737.75
$ws.Range("J100").Value = 1500.1666
$ws.Range("K100").Value = 1737.75
$ws.Range("L100").Value = 1500.1666
$ws.Range("M100").Value = -1196.75
$ws.Range("N100").Value = -2582.1666
$ws.Range("H132").Value = 402583.25
$ws.Range("I132").Value = 2593.4443
$ws.Range("J132").Value = 1431128.4
$ws.Range("K132").Value = 7780.3329
$ws.Range("L132").Value = 4293385.199999999
$ws.Range("M132").Value = -5250.3329
$ws.Range("N132").Value = -4298445.199999999
$ws.Range("H137").Value = 1019.7692
$ws.Range("I137").Value = 868.5455
$ws.Range("J137").Value = 1851.5
$ws.Range("K137").Value = 2605.6365
$ws.Range("L137").Value = 5554.5
$ws.Range("M137").Value = -55.63649999999961
$ws.Range("N137").Value = -10654.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 19000
$ws.Range("J9").Value = 19000
$ws.Range("L9").Value = 19000
$ws.Range("N9").Value = -19340
$ws.Range("H20").Value = 19000
$ws.Range("J20").Value = 19000
$ws.Range("L20").Value = 19000
$ws.Range("N20").Value = -19540
$ws.Range("H45").Value = 2185.3333
$ws.Range("I45").Value = 1913.875
$ws.Range("J45").Value = 2728.25
$ws.Range("K45").Value = 1913.875
$ws.Range("L45").Value = 2728.25
$ws.Range("M45").Value = -1536.875
$ws.Range("N45").Value = -3482.25
$ws.Range("H88").Value = 3160
$ws.Range("I88").Value = 2700
$ws.Range("J88").Value = 3850
$ws.Range("K88").Value = 2700
$ws.Range("L88").Value = 3850
$ws.Range("M88").Value = -2294
$ws.Range("N88").Value = -4662
$ws.Range("H91").Value = 3160
$ws.Range("I91").Value = 2700
$ws.Range("J91").Value = 3850
$ws.Range("K91").Value = 2700
$ws.Range("L91").Value = 3850
$ws.Range("M91").Value = -1296
$ws.Range("N91").Value = -6658
$ws.Range("H139").Value = 49813.332
$ws.Range("J139").Value = 49813.332
$ws.Range("L139").Value = 49813.332
$ws.Range("N139").Value = -60093.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2324.375
$ws.Range("I86").Value = 2208.182
$ws.Range("J86").Value = 2580
$ws.Range("K86").Value = 2208.182
$ws.Range("L86").Value = 2580
$ws.Range("M86").Value = -1085.182
$ws.Range("N86").Value = -4826
$ws.Range("H89").Value = 2324.375
$ws.Range("I89").Value = 2208.182
$ws.Range("J89").Value = 2580
$ws.Range("K89").Value = 11040.91
$ws.Range("L89").Value = 12900
$ws.Range("M89").Value = -5424.91
$ws.Range("N89").Value = -24132
$ws.Range("H137").Value = 55338.46
$ws.Range("J137").Value = 55338.46
$ws.Range("L137").Value = 55338.46
$ws.Range("N137").Value = -65538.45999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2958.318
$ws.Range("I31").Value = 2958.318
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2958.318
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2663.318
$ws.Range("N31").Value = $null
$ws.Range("H34").Value = 2958.318
$ws.Range("I34").Value = 2958.318
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2958.318
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2756.318
$ws.Range("N34").Value = $null
$ws.Range("H62").Value = 9846.471
$ws.Range("I62").Value = 11550
$ws.Range("J62").Value = 8332.223
$ws.Range("K62").Value = 11550
$ws.Range("L62").Value = 8332.223
$ws.Range("M62").Value = -10926
$ws.Range("N62").Value = -9580.223
$ws.Range("H65").Value = 9846.471
$ws.Range("I65").Value = 11550
$ws.Range("J65").Value = 8332.223
$ws.Range("K65").Value = 57750
$ws.Range("L65").Value = 41661.115
$ws.Range("M65").Value = -54630
$ws.Range("N65").Value = -47901.115
$ws.Range("H105").Value = 2669.6667
$ws.Range("I105").Value = 3340
$ws.Range("J105").Value = 1999.3334
$ws.Range("K105").Value = 3340
$ws.Range("L105").Value = 1999.3334
$ws.Range("M105").Value = -1593
$ws.Range("N105").Value = -5493.3334
$ws.Range("H107").Value = 365.86365
$ws.Range("I107").Value = 234
$ws.Range("J107").Value = 404.64706
$ws.Range("K107").Value = 234
$ws.Range("L107").Value = 404.64706
$ws.Range("M107").Value = 1686
$ws.Range("N107").Value = -4244.64706

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 648714.7
$ws.Range("I107").Value = 529.9
$ws.Range("K107").Value = 1589.7
$ws.Range("M107").Value = 330.3000000000002
$ws.Range("H117").Value = 2520
$ws.Range("I117").Value = 964.5
$ws.Range("J117").Value = 2759.3076
$ws.Range("K117").Value = 2893.5
$ws.Range("L117").Value = 8277.9228
$ws.Range("M117").Value = 548.5
$ws.Range("N117").Value = -15161.9228
$ws.Range("H122").Value = 556546.8
$ws.Range("I122").Value = 576
$ws.Range("J122").Value = 1112517.6
$ws.Range("K122").Value = 5184
$ws.Range("L122").Value = 10012658.4
$ws.Range("M122").Value = -2734
$ws.Range("N122").Value = -10017558.4
$ws.Range("H129").Value = 2399.25
$ws.Range("I129").Value = 1380
$ws.Range("J129").Value = 2965.5
$ws.Range("K129").Value = 4140
$ws.Range("L129").Value = 8896.5
$ws.Range("M129").Value = 860
$ws.Range("N129").Value = -18896.5
$ws.Range("H131").Value = 829.74225
$ws.Range("I131").Value = 431.26666
$ws.Range("J131").Value = 902.63416
$ws.Range("K131").Value = 1293.79998
$ws.Range("L131").Value = 2707.90248
$ws.Range("M131").Value = 3746.20002
$ws.Range("N131").Value = -12787.90248

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4375
$ws.Range("I126").Value = 4375
$ws.Range("K126").Value = 13125
$ws.Range("M126").Value = -10655

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1246.6666
$ws.Range("I46").Value = 1494.2858
$ws.Range("K46").Value = 1494.2858
$ws.Range("M46").Value = -1306.2858
$ws.Range("H132").Value = 4660.1333
$ws.Range("I132").Value = 4863
$ws.Range("K132").Value = 14589
$ws.Range("M132").Value = -12059

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 7753.1
$ws.Range("J41").Value = 7753.1
$ws.Range("L41").Value = 7753.1
$ws.Range("N41").Value = -8533.1
